# 2025 roswell data added
# Add a new drug name entry ("trulicity(4)") as the next row after the
# existing list (A1:A26 -> A1:A27), then move the selection to the next
# empty row (A28), matching where Excel would leave the cursor after
# typing a new entry and pressing Enter.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A27").Value = "trulicity(4)"

$ws.Range("A28").Select()
